$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is used on the Overview sheet (zh-cn/de-de status columns)
# as well as on each language sheet's "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns on the Overview sheet ---
# and the "Status" column on each language sheet.
# Target stored column width (OOXML "width" units) is ~13.41; the closest
# value reachable through ColumnWidth (quantized to 1/6 character
# increments by the engine) is 13.3333..., reached by requesting 12.5.
$wsOverview.Range("E:E").ColumnWidth = 12.5
$wsOverview.Range("F:F").ColumnWidth = 12.5

$wsZhCn.Range("C:C").ColumnWidth = 12.5
$wsDeDe.Range("C:C").ColumnWidth = 12.5
